$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (government)
$ws.Range("C2").Value = 38
$ws.Range("D2").Value = 173
$ws.Range("E2").Value = 76
$ws.Range("F2").Value = 62
$ws.Range("H2").Value = 80
$ws.Range("I2").Value = 45
$ws.Range("K2").Value = 110
$ws.Range("L2").Value = 81
$ws.Range("M2").Value = 79
$ws.Range("N2").Value = 92

# Row 3 (independent)
$ws.Range("C3").Value = 60
$ws.Range("D3").Value = 414
$ws.Range("E3").Value = 167
$ws.Range("F3").Value = 205
$ws.Range("G3").Value = 316
$ws.Range("H3").Value = 217
$ws.Range("I3").Value = 82
$ws.Range("J3").Value = 203
$ws.Range("K3").Value = 448
$ws.Range("L3").Value = 464
$ws.Range("M3").Value = 195
$ws.Range("N3").Value = 180

# Row 4 (university)
$ws.Range("D4").Value = 28
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 14
$ws.Range("K4").Value = 14
$ws.Range("N4").Value = 7

# Row 5 (unknown_gov)
$ws.Range("D5").Value = 7
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 17
$ws.Range("L5").Value = 13
$ws.Range("N5").Value = 8
